$wb = $excel.ActiveWorkbook

# --- Summary sheet: bump fee figures from 50 to 100 and move the selection ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("A4").Value = 100
$wsSummary.Range("B4").Value = 100
[void]$wsSummary.Range("B4").Select()

# --- Repayment Schedule sheet: bump fee figures from 50 to 100 and move the selection ---
$wsRepay = $wb.Worksheets.Item("Repayment Schedule")
$wsRepay.Range("I2").Value = 100
$wsRepay.Range("K2").Value = 100
$wsRepay.Range("L2").Value = 100
[void]$wsRepay.Range("L2").Select()

# --- Transactions sheet: bump fee figures from 50 to 100, make this the active tab ---
$wsTrans = $wb.Worksheets.Item("Transactions")
$wsTrans.Range("E2").Value = 100
$wsTrans.Range("H2").Value = 100
[void]$wsTrans.Activate()
[void]$wsTrans.Range("I13:J13").Select()
